$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4797
$ws.Range("I32").Value = 5749.25
$ws.Range("J32").Value = 4162.1665
$ws.Range("K32").Value = 5749.25
$ws.Range("L32").Value = 4162.1665
$ws.Range("M32").Value = -5423.25
$ws.Range("N32").Value = -4814.1665
$ws.Range("H80").Value = 546.9286
$ws.Range("I80").Value = 374.8889
$ws.Range("J80").Value = 856.6
$ws.Range("K80").Value = 1124.6667
$ws.Range("L80").Value = 2569.8
$ws.Range("M80").Value = -126.6667
$ws.Range("N80").Value = -4565.8
$ws.Range("H83").Value = 546.9286
$ws.Range("I83").Value = 374.8889
$ws.Range("J83").Value = 856.6
$ws.Range("K83").Value = 3374.0001
$ws.Range("L83").Value = 7709.400000000001
$ws.Range("M83").Value = 1617.9999
$ws.Range("N83").Value = -17693.4
$ws.Range("H96").Value = 1459.4546
$ws.Range("I96").Value = 1756.875
$ws.Range("J96").Value = 666.3333
$ws.Range("K96").Value = 5270.625
$ws.Range("L96").Value = 1998.9999
$ws.Range("M96").Value = -3897.625
$ws.Range("N96").Value = -4744.9999
$ws.Range("H104").Value = 116.083336
$ws.Range("I104").Value = 115.72727
$ws.Range("K104").Value = 347.18181
$ws.Range("M104").Value = 1399.81819
$ws.Range("H135").Value = 551.0454999999999
$ws.Range("I135").Value = 539.35
$ws.Range("K135").Value = 4854.150000000001
$ws.Range("M135").Value = -2319.150000000001
$ws.Range("H138").Value = 4622.963
$ws.Range("J138").Value = 6082
$ws.Range("L138").Value = 18246
$ws.Range("N138").Value = -28526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6358.615
$ws.Range("I61").Value = 7860.6113
$ws.Range("K61").Value = 7860.6113
$ws.Range("M61").Value = -7648.6113
$ws.Range("H92").Value = 74772
$ws.Range("J92").Value = 74772
$ws.Range("L92").Value = 74772
$ws.Range("N92").Value = -79764
$ws.Range("H102").Value = 4389686.5
$ws.Range("I102").Value = 4389686.5
$ws.Range("K102").Value = 4389686.5
$ws.Range("M102").Value = -4388064.5
$ws.Range("H136").Value = 6358.615
$ws.Range("I136").Value = 7860.6113
$ws.Range("K136").Value = 23581.8339
$ws.Range("M136").Value = -21031.8339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3956707.2
$ws.Range("I94").Value = 6993403.5
$ws.Range("K94").Value = 6993403.5
$ws.Range("M94").Value = -6992952.5
$ws.Range("H105").Value = 4234446.5
$ws.Range("I105").Value = 4234446.5
$ws.Range("K105").Value = 4234446.5
$ws.Range("M105").Value = -4232699.5
$ws.Range("H134").Value = 9969.540999999999
$ws.Range("I134").Value = 9968.107
$ws.Range("K134").Value = 29904.321
$ws.Range("M134").Value = -27369.321

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 33666.668
$ws.Range("I44").Value = 90000
$ws.Range("K44").Value = 90000
$ws.Range("M44").Value = -89558
$ws.Range("H60").Value = 33333
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 44999.5
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 44999.5
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -46021.5
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 21456.182
$ws.Range("I87").Value = 20006
$ws.Range("K87").Value = 60018
$ws.Range("M87").Value = -58770
$ws.Range("H90").Value = 21456.182
$ws.Range("I90").Value = 20006
$ws.Range("K90").Value = 180054
$ws.Range("M90").Value = -173814
$ws.Range("H121").Value = 249.57143
$ws.Range("I121").Value = 249.57143
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 748.71429
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 561.28571
$ws.Range("N121").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H132").Value = 881.7692
$ws.Range("I132").Value = 846.9
$ws.Range("J132").Value = 998
$ws.Range("K132").Value = 7622.099999999999
$ws.Range("L132").Value = 8982
$ws.Range("M132").Value = -5092.099999999999
$ws.Range("N132").Value = -14042

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 48258.25
$ws.Range("J52").Value = 48258.25
$ws.Range("L52").Value = 48258.25
$ws.Range("N52").Value = -48776.25
$ws.Range("H80").Value = 1442251.9
$ws.Range("I80").Value = 2042426
$ws.Range("J80").Value = 1834
$ws.Range("K80").Value = 2042426
$ws.Range("L80").Value = 1834
$ws.Range("M80").Value = -2041428
$ws.Range("N80").Value = -3830
$ws.Range("H83").Value = 1442251.9
$ws.Range("I83").Value = 2042426
$ws.Range("J83").Value = 1834
$ws.Range("K83").Value = 10212130
$ws.Range("L83").Value = 9170
$ws.Range("M83").Value = -10207138
$ws.Range("N83").Value = -19154
$ws.Range("H97").Value = 993167.5600000001
$ws.Range("I97").Value = 1701817.8
$ws.Range("J97").Value = 1057.3
$ws.Range("K97").Value = 1701817.8
$ws.Range("L97").Value = 1057.3
$ws.Range("M97").Value = -1701321.8
$ws.Range("N97").Value = -2049.3
$ws.Range("H102").Value = 9064921
$ws.Range("I102").Value = 27779526
$ws.Range("J102").Value = 2259610
$ws.Range("K102").Value = 27779526
$ws.Range("L102").Value = 2259610
$ws.Range("M102").Value = -27777904
$ws.Range("N102").Value = -2262854
$ws.Range("H122").Value = 685326.25
$ws.Range("H132").Value = 9379.031000000001
$ws.Range("I132").Value = 6643.24
$ws.Range("K132").Value = 19929.72
$ws.Range("M132").Value = -17399.72

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8581428
$ws.Range("I2").Value = 15002500
$ws.Range("J2").Value = 19999.666
$ws.Range("K2").Value = 15002500
$ws.Range("L2").Value = 19999.666
$ws.Range("M2").Value = -15002388
$ws.Range("N2").Value = -20223.666
$ws.Range("H10").Value = 84499.586
$ws.Range("I10").Value = 2003
$ws.Range("J10").Value = 199994.8
$ws.Range("K10").Value = 2003
$ws.Range("L10").Value = 199994.8
$ws.Range("M10").Value = -1863
$ws.Range("N10").Value = -200274.8
$ws.Range("H45").Value = 16688.5
$ws.Range("I45").Value = 17729.875
$ws.Range("J45").Value = 12523
$ws.Range("K45").Value = 17729.875
$ws.Range("L45").Value = 12523
$ws.Range("M45").Value = -17322.875
$ws.Range("N45").Value = -13337
$ws.Range("H48").Value = 25000
$ws.Range("H55").Value = 1189.2307
$ws.Range("I55").Value = 921.25
$ws.Range("K55").Value = 921.25
$ws.Range("M55").Value = -748.25
$ws.Range("H82").Value = 4275491
$ws.Range("I82").Value = 6174822
$ws.Range("J82").Value = 1996.5
$ws.Range("K82").Value = 6174822
$ws.Range("L82").Value = 1996.5
$ws.Range("M82").Value = -6174461
$ws.Range("N82").Value = -2718.5
$ws.Range("H85").Value = 4275491
$ws.Range("I85").Value = 6174822
$ws.Range("J85").Value = 1996.5
$ws.Range("K85").Value = 6174822
$ws.Range("L85").Value = 1996.5
$ws.Range("M85").Value = -6173574
$ws.Range("N85").Value = -4492.5
$ws.Range("H122").Value = 7112.2
$ws.Range("I122").Value = 3986.25
$ws.Range("K122").Value = 11958.75
$ws.Range("M122").Value = -9508.75
$ws.Range("H132").Value = 11954.619
$ws.Range("J132").Value = 9436.75
$ws.Range("L132").Value = 28310.25
$ws.Range("N132").Value = -33370.25
$ws.Range("H136").Value = 42231.035
$ws.Range("I136").Value = 54212.95
$ws.Range("K136").Value = 162638.85
$ws.Range("M136").Value = -160088.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 65000
$ws.Range("J70").Value = 65000
$ws.Range("L70").Value = 65000
$ws.Range("N70").Value = -65630
$ws.Range("H73").Value = 65000
$ws.Range("J73").Value = 65000
$ws.Range("L73").Value = 65000
$ws.Range("N73").Value = -67184
$ws.Range("H75").Value = 40065
$ws.Range("J75").Value = 40065
$ws.Range("L75").Value = 40065
$ws.Range("N75").Value = -41937
$ws.Range("H78").Value = 40065
$ws.Range("J78").Value = 40065
$ws.Range("L78").Value = 120195
$ws.Range("N78").Value = -129555
$ws.Range("H96").Value = 2985.5715
$ws.Range("J96").Value = 2985.5715
$ws.Range("L96").Value = 2985.5715
$ws.Range("N96").Value = -5731.5715
$ws.Range("H100").Value = 830.8889
$ws.Range("I100").Value = 166.33333
$ws.Range("K100").Value = 332.66666
$ws.Range("M100").Value = 208.33334
$ws.Range("H113").Value = 869.25
$ws.Range("J113").Value = 1088.5
$ws.Range("L113").Value = 3265.5
$ws.Range("N113").Value = -7605.5
$ws.Range("H122").Value = 3257.647
$ws.Range("I122").Value = 3189
$ws.Range("J122").Value = 3355.7144
$ws.Range("K122").Value = 9567
$ws.Range("L122").Value = 10067.1432
$ws.Range("M122").Value = -7117
$ws.Range("N122").Value = -14967.1432
$ws.Range("H136").Value = 7000.816
$ws.Range("J136").Value = 6778.75
$ws.Range("L136").Value = 20336.25
$ws.Range("N136").Value = -25436.25
